# ENA - Raw sequencing reads: update version + Term Source/Accession refs (TSR update)
$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump Version 1.0.2 -> 1.0.3 ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.3"

# --- New Table sheet: table header + data updates ---
$wsTable = $wb.Worksheets.Item("New Table")

# Header row (also drives the annotationTable column names in table1.xml)
$wsTable.Range("B1").Value = "Component [next generation sequencing instrument model]"
$wsTable.Range("C1").Value = "Term Source REF (DPBO:0000040)"
$wsTable.Range("D1").Value = "Term Accession Number (DPBO:0000040)"

# Data row: Instrument Model ontology now GENEPIO instead of OBI
$wsTable.Range("C2").Value = "GENEPIO"
$wsTable.Range("D2").Value = "http://purl.obolibrary.org/obo/GENEPIO_0100115"

# Data row: Library selection method renamed, now sourced from GENEPIO
$wsTable.Range("K2").Value = "PCR method"
$wsTable.Range("L2").Value = "GENEPIO"
$wsTable.Range("M2").Value = "http://purl.obolibrary.org/obo/GENEPIO_0001955"
